$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11 ("Tues, Feb 23" / CI: Matching lecture) ---
# "In-Class Exercise" (D11): drop the Slides bullet, keep only the Exercises link.
$ws.Range("D11").Value = '- `Exercises <exercises/exercise_matching.ipynb>`_'

# "Do Before Class" (C11): add a new DAME-FLAME User Guide reading alongside
# the existing Matching Methods reading.
$c11_1 = '- `Matching Methods for Causal Inference: A Review and a Look Forward <https://www.ncbi.nlm.nih.gov/pmc/articles/PMC2943670/>`_'
$c11_2 = '- `DAME-FLAME User Guide <https://almost-matching-exactly.github.io/DAME-FLAME-Python-Package/user-guide>`_'
$ws.Range("C11").Value = $c11_1 + "`n" + $c11_2

# --- Row 6 ("Thurs, Feb 4" / Internal v. External Validity lecture) ---
# "Do Before Class" (C6): turn the plain heading into a link, rename the
# "Endogenous Stopping Times" reading, and add a new A/B Testing reading.
$c6_1 = '- `Internal versus External Validity <internal_v_external_validity.ipynb>`_'
$c6_2 = '- `Limitations of Average Treatment Effects <limitations_of_ATE.ipynb>`_'
$c6_3 = '- `A/B Testing and Stopping Early <https://medium.com/airbnb-engineering/experiments-at-airbnb-e2db3abf39e7>`_'
$c6_4 = '- Optional: `More on endogenous stopping <http://sl8r000.github.io/ab_testing_statistics/avoid_biased_stopping_times/>`_'
$c6_5 = '- Optional: `Discussion in part About Experiments in Advertising <https://overcast.fm/+QLduPjO1I>`_'
$ws.Range("C6").Value = $c6_1 + "`n" + $c6_2 + "`n" + $c6_3 + "`n" + $c6_4 + "`n" + $c6_5
